$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Абраменкова
$ws.Cells.Item(2,2).Value = 4

# Row 3 - Агапова
$ws.Cells.Item(3,2).Value = -2
$ws.Cells.Item(3,4).Value = 0

# Row 4 - Бычкова
$ws.Cells.Item(4,2).Value = 4
$ws.Cells.Item(4,3).Value = 4
$ws.Cells.Item(4,13).Value = "переписаны верно все номера"

# Row 5 - Виноградов
$ws.Cells.Item(5,3).Value = 4
$ws.Cells.Item(5,4).Value = 5
$ws.Cells.Item(5,13).Value = "переписаны верно все номера"

# Row 6 - Воробьев
$ws.Cells.Item(6,4).Value = 5

# Row 7 - Глазков
$ws.Cells.Item(7,4).Value = 5

# Row 8 - Губеева
$ws.Cells.Item(8,2).Value = 4
$ws.Cells.Item(8,4).Value = 5

# Row 9 - Заднипрянец
$ws.Cells.Item(9,2).Value = 4

# Row 10 - Ибрамхалилов
$ws.Cells.Item(10,4).Value = 5

# Row 11 - Капелина
$ws.Cells.Item(11,4).Value = 5

# Row 13 - Куулар
$ws.Cells.Item(13,2).Value = 4
$ws.Cells.Item(13,3).Value = 4
$ws.Cells.Item(13,4).Value = -1
$ws.Cells.Item(13,13).Value = "переписаны верно все номера"

# Row 15 - Лунин
$ws.Cells.Item(15,5).Value = 5
$ws.Cells.Item(15,7).Value = 0

# Row 16 - Оганезов
$ws.Cells.Item(16,4).Value = 5

# Row 17 - Попова
$ws.Cells.Item(17,3).Value = 4
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,13).Value = "переписаны верно все номера"

# Row 18 - Родина
$ws.Cells.Item(18,4).Value = 5

# Row 19 - Словогородская
$ws.Cells.Item(19,4).Value = 5

# Row 20 - Стоценко
$ws.Cells.Item(20,4).Value = 5

# Row 22 - Титова
$ws.Cells.Item(22,2).Value = 4
$ws.Cells.Item(22,4).Value = 5

# Row 23 - Хрищанович
$ws.Cells.Item(23,2).Value = 4
$ws.Cells.Item(23,4).Value = 5

# Row 24 - Чиченкова
$ws.Cells.Item(24,4).Value = 5

# Update the active cell selection to D4
$ws.Range("D4").Select()
